# Append: 2025-10-02 06:25 JST
# Replaces the "ランサーズ" worksheet's data rows (2-21) with a fresh set of
# 5 listings (rows 2-6), shrinks a few column widths, and rewrites the
# F-column hyperlinks so they point at the new listing URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the old rows (7-21) that no longer exist in the refreshed data.
#    This also shrinks UsedRange/dimension down to A1:H6 automatically.
# ---------------------------------------------------------------------
$ws.Rows("7:21").Delete()

# ---------------------------------------------------------------------
# 2. Clear out every existing hyperlink (stale rIds for the removed rows
#    otherwise linger in the package) - we will add back fresh ones below.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3. New row data (row 1 header is untouched).
# ---------------------------------------------------------------------
$timestamp = "2025-10-02 06:25:32"

$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5405023"
$ws.Range("G2").Value = 178
$ws.Range("H2").Value = "★bot ◆ツール"

$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】データ処理のためのExcel VBA・マクロ開発依頼 もしくはスクレイピングによる対応"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5405218"
$ws.Range("G3").Value = 98
$ws.Range("H3").Value = "◆開発,スクレイピング"

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【急募】集計分析ツール(keyence社製「KI」)の設定構築経験者"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5405052"
$ws.Range("G4").Value = 73
$ws.Range("H4").Value = "◆ツール"

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【在宅勤務】ランサーズ業務委託で働ける、ネパール人個人の方を募集します!"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5404906"
$ws.Range("G5").Value = 18
$ws.Range("H5").ClearContents()

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "LINE公式アカウントの動作確認・タグ等設定対応"
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5405235"
$ws.Range("G6").Value = 10
$ws.Range("H6").ClearContents()

# ---------------------------------------------------------------------
# 4. Re-create the hyperlinks for the URL column (F2:F6), in order, so
#    the relationship ids come out as rId1..rId5 pointing at the right
#    targets.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405023")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5405218")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5405052")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5404906")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5405235")

# ---------------------------------------------------------------------
# 5. Column width tweaks (B 52->49, D 32->27, H 23->13).
#    Excel's ColumnWidth setter stores width + 5/6 in the XML, so we
#    subtract that constant offset up front to land on the exact target.
# ---------------------------------------------------------------------
$offset = 5.0 / 6.0
$ws.Columns.Item(2).ColumnWidth = 49 - $offset
$ws.Columns.Item(4).ColumnWidth = 27 - $offset
$ws.Columns.Item(8).ColumnWidth = 13 - $offset
